# spring 24 week 5 inputs
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = 1.24
$ws.Range("E3").Value = 1.32
$ws.Range("F3").Value = 1.23
$ws.Range("B4").Value = 1.49
$ws.Range("C4").Value = 1.45
$ws.Range("F4").Value = 1.11
$ws.Range("C5").Value = 1.33
$ws.Range("F5").Value = 1.04
$ws.Range("C6").Value = 1.43
$ws.Range("D6").Value = 1.49
$ws.Range("E6").Value = 1.32
$ws.Range("D7").Value = 1.77
